$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns I and J, formatted like the existing header cells
$ws.Range("I1").Value = "I0"
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)

$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Data rows
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 8

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 4

$ws.Range("I4").Value = 7
$ws.Range("J4").Value = 8

$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 3
